# Commit: "fixed auto, added trainer in gui"
# Adds a new trainer row ("soham") to the details sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (Sr. No, Name, Address, Job) below the existing data.
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "soham"
$ws.Range("C6").Value = "202/Ring city road, Bandra,Mumbai"
$ws.Range("D6").Value = "Deployment"

# Update the active selection to match the post-edit cursor position.
$ws.Range("E9").Select()
